$d = $word.ActiveDocument

# --- Part 1: "AE/VSA " paragraph -------------------------------------------
# Remove "- 2/14  Some VSA Contract recompete details" (plain text + hyperlink
# run) and trim the trailing space off "AE/VSA ".
$rng = $d.Content
$rng.Find.Execute("- 2/14  Some VSA Contract recompete details", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

$rng2 = $d.Content
$rng2.Find.Execute("AE/VSA ", $false, $false, $false, $false, $false, $true, 1, $false, "AE/VSA", 2)

# --- Part 2: collapse the "Transition" bullet list --------------------------
# Find the paragraph that now reads "Transition documents had to be created..."
# (originally the 3rd bullet in this run, right after the image paragraph) and
# the paragraph that follows it, through the blank line right before the
# "#31611 My VA Enhancements" heading. Delete the image/caption paragraph that
# precedes it, and the five bullet paragraphs (plus blank line) that follow
# it, leaving only the "Transition documents..." bullet behind.

function Find-ParaIndex($text) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $t = $p.Range.Text.TrimEnd("`r")
        if ($t -eq $text) {
            return $i
        }
    }
    return $null
}

$target = Find-ParaIndex("Transition documents had to be created, even though we were staying - Links at the bottom of this document")
Write-Output "target1=$target"

# Delete the paragraph right before it (the "Auth Exp..." + image paragraph).
$paras = $d.Paragraphs
$before = $paras.Item($target - 1)
$before.Range.Delete()

# Re-fetch paragraphs/index since the collection shifted after the delete.
$target = Find-ParaIndex("Transition documents had to be created, even though we were staying - Links at the bottom of this document")
Write-Output "target2=$target"

# Delete the five paragraphs after it (Prime/4-1/Adhoc/Org Changes/Transition
# Timeline) plus the trailing blank paragraph, all the way up to (but not
# including) the following heading paragraph.
$paras = $d.Paragraphs
$afterStart = $paras.Item($target + 1)
$afterEnd = $paras.Item($target + 6)
$delRange = $d.Range($afterStart.Range.Start, $afterEnd.Range.End)
Write-Output "delText=[$($delRange.Text)]"
$delRange.Delete()

# Update the spacing-after of the remaining "Transition documents..." bullet
# to match what the last item in the list used to have (240 twips -> 12pt).
$target = Find-ParaIndex("Transition documents had to be created, even though we were staying - Links at the bottom of this document")
Write-Output "target3=$target"
$paras = $d.Paragraphs
$p = $paras.Item($target)
$p.SpaceAfter = 12

Write-Output "Done"
